# Weekly update: insert a new price-report row for "Bruselas (repollito)"
# at the top of the data block (row 14), pushing the existing rows 14-71
# down to 15-72. The new row carries the same constant categorical fields
# as the rest of the block, with a fresh date and updated price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 14; everything below shifts down by one.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Range("A14").Value = 6
$ws.Range("B14").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = 45071
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112035
$ws.Range("G14").Value = "Bruselas (repollito)"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 290
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 22000
$ws.Range("M14").Value = 21172
$ws.Range("N14").Value = "$/malla 15 kilos"
$ws.Range("O14").Value = "Provincia de Quillota"
$ws.Range("P14").Value = 1411
$ws.Range("Q14").Value = 15
$ws.Range("R14").Value = "Hortaliza"
